# Implement type designator using 'type' field for technique-specific sample
# preparation sheets. Adds the SamplePreparation "header" columns (type,
# sample_id, preparation_date, operator_id, protocol_description, id, title)
# in front of the existing trailing "description" column on each of the
# technique-specific preparation sheets, and renames SamplePreparation's own
# "preparation_type" column to "type" while dropping its dropdown validation
# (the enum is no longer needed now that 'type' is a free-form designator).

$wb = $excel.ActiveWorkbook

# Columns shared by the technique-specific preparation sheets, inserted right
# before the trailing "description" column.
$newColumns = @("type", "sample_id", "preparation_date", "operator_id", "protocol_description", "id", "title")

# sheetName -> column letter currently holding "description" (i.e. the last
# used column before this edit).
$targets = @{
    "CryoEMPreparation" = "J"
    "XRayPreparation"   = "H"
    "SAXSPreparation"   = "F"
}

foreach ($sheetName in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $descCol = $targets[$sheetName]

    # Remember the current description text (normally "description"),
    # then overwrite that cell with the first new column and append the
    # rest, finishing with description restored one column further right.
    $descValue = $ws.Range($descCol + "1").Value2

    $col = $ws.Range($descCol + "1").Column
    foreach ($name in $newColumns) {
        $ws.Cells.Item(1, $col).Value = $name
        $col = $col + 1
    }
    $ws.Cells.Item(1, $col).Value = $descValue
}

# SamplePreparation: rename preparation_type -> type, drop its A-column list
# validation (cryo_em, xray_crystallography, ... enum goes away in favor of
# free-form technique type designators).
$sp = $wb.Worksheets.Item("SamplePreparation")
$sp.Range("A1").Value = "type"
$sp.Columns.Item(1).Validation.Delete()
